$wb = $excel.ActiveWorkbook
$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# --- Timestamp update (column AK, rows 2-48, sheet FBS) ---
for ($r = 2; $r -le 48; $r++) {
    $wsFBS.Cells.Item($r, 37).Value = "2025-09-26T10:01:49.321475"
}

# --- Per-cell data updates ---
$wsFBS.Range("Q4").Value = "WNW"
$wsFBS.Range("Z4").Value = -115
$wsFBS.Range("N10").Value = "SSW"
$wsFBS.Range("Q11").Value = "E"
$wsFBS.Range("Y11").Value = 63.5
$wsFBS.Range("Z11").Value = -114
$wsFBS.Range("AE11").Value = 0.03252032520325204
$wsFBS.Range("Q15").Value = "E"
$wsFBS.Range("Y15").Value = 60.5
$wsFBS.Range("Z15").Value = -110
$wsFBS.Range("AE15").Value = 0.03418803418803419
$wsFBS.Range("N16").Value = "ESE"
$wsFBS.Range("N17").Value = "SSW"
$wsFBS.Range("O17").Value = 76.52
$wsFBS.Range("P17").Value = 4
$wsFBS.Range("U17").Value = 0.3
$wsFBS.Range("N18").Value = "WSW"
$wsFBS.Range("O18").Value = 71.48
$wsFBS.Range("P18").Value = 11.4
$wsFBS.Range("Q18").Value = "WSW"
$wsFBS.Range("R18").Value = 5.499999999999999
$wsFBS.Range("U18").Value = 5.1
$wsFBS.Range("M19").Value = "SSW"
$wsFBS.Range("O19").Value = 63.88999999999999
$wsFBS.Range("P19").Value = 5.2
$wsFBS.Range("U19").Value = 0.2
$wsFBS.Range("M20").Value = "SSW"
$wsFBS.Range("M21").Value = "ESE"
$wsFBS.Range("N22").Value = "W"
$wsFBS.Range("Z22").Value = -114
$wsFBS.Range("N23").Value = "SSW"
$wsFBS.Range("Q23").Value = "SSW"
$wsFBS.Range("Q25").Value = "E"
$wsFBS.Range("Z26").Value = -112
$wsFBS.Range("N28").Value = "SSW"
$wsFBS.Range("Q28").Value = "SSW"
$wsFBS.Range("R28").Value = 0.7
$wsFBS.Range("M29").Value = "W"
$wsFBS.Range("R29").Value = 0
$wsFBS.Range("Q31").Value = "NE"
$wsFBS.Range("Q33").Value = "NE"
$wsFBS.Range("Q34").Value = "SE"
$wsFBS.Range("N35").Value = "W"
$wsFBS.Range("Q35").Value = "W"
$wsFBS.Range("R35").Value = 5.8
$wsFBS.Range("M36").Value = "SSW"
$wsFBS.Range("Q36").Value = "SSW"
$wsFBS.Range("Z36").Value = -110
$wsFBS.Range("R37").Value = 20.2
$wsFBS.Range("Z37").Value = -118
$wsFBS.Range("N40").Value = "W"
$wsFBS.Range("Q40").Value = "W"
$wsFBS.Range("Z40").Value = -105
$wsFBS.Range("M42").Value = "ESE"
$wsFBS.Range("Q42").Value = "E"
$wsFBS.Range("Q43").Value = "E"
$wsFBS.Range("M44").Value = "W"
$wsFBS.Range("Q44").Value = "W"
$wsFBS.Range("N47").Value = "SSW"
$wsOther.Range("P3").Value = "ESE"
$wsOther.Range("S3").Value = "E"
$wsOther.Range("O5").Value = "W"
$wsOther.Range("P5").Value = "W"
$wsOther.Range("S5").Value = "W"
$wsOther.Range("T8").Value = 0.7
$wsOther.Range("P10").Value = "SSW"
$wsOther.Range("S10").Value = "SSW"
$wsOther.Range("S12").Value = "NNE"
$wsOther.Range("P13").Value = "W"
$wsOther.Range("O14").Value = "W"
$wsOther.Range("S14").Value = "W"
$wsOther.Range("O15").Value = "W"
$wsOther.Range("O17").Value = "W"
$wsOther.Range("S17").Value = "W"
$wsOther.Range("O20").Value = "W"
$wsOther.Range("S20").Value = "W"
$wsOther.Range("S22").Value = "NNW"
$wsOther.Range("O24").Value = "SSW"
$wsOther.Range("P24").Value = "SSW"
$wsOther.Range("S24").Value = "SSW"
$wsOther.Range("T24").Value = 22.3
$wsOther.Range("T32").Value = 0
$wsOther.Range("O34").Value = "ESE"
$wsOther.Range("S34").Value = "E"
$wsOther.Range("T34").Value = 61.5
$wsOther.Range("T35").Value = 0.7
$wsOther.Range("S37").Value = "E"
$wsOther.Range("T37").Value = 0.4
$wsOther.Range("P39").Value = "SSW"
$wsOther.Range("S39").Value = "SSW"
$wsOther.Range("S42").Value = "E"
$wsOther.Range("S44").Value = "NNE"
